# Insert a new event row ("Note di Stelle") at row 27 of the "eventi" sheet.
# This shifts all existing rows from 27..55 down to 28..56, and the new
# row 27 is populated with the new event's data (mirroring the previous
# row 27's layout/empty-cell pattern, since no style/row attributes are
# otherwise involved in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 27 (and everything below it) down by one row.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the "Note di Stelle" event data.
$ws.Range("A27").Value2 = "Altri eventi,Spettacoli,Musica"
$ws.Range("B27").Value2 = "Modena"
$ws.Range("C27").Value2 = "Centro storico"
$ws.Range("D27").Value2 = "2022-06-03T15:42:40+00:00"
$ws.Range("E27").Value2 = "Rassegna di spettacoli in Piazza XX Settembre"
$ws.Range("F27").Value2 = "2022-06-03T15:43:48+00:00"
$ws.Range("G27").Value2 = ""
$ws.Range("H27").Value2 = "2022-06-07T15:00:00+00:00"
$ws.Range("I27").Value2 = "2022-06-28T16:00:00+00:00"
$ws.Range("J27").Value2 = "https://www.comune.modena.it/api/novita/eventi/2022/note-di-stelle/@@images/78120e75-8036-4a4c-b46e-10cc1ff7fe4c.jpeg"
$ws.Range("K27").Value2 = ""
$ws.Range("L27").Value2 = "2022-06-03T15:43:48+00:00"
$ws.Range("M27").Value2 = "Piazza XX Settembre"
$ws.Range("N27").Value2 = " 7, 14, 21 e 28 giugno alle ore 21.00"
$ws.Range("O27").Value2 = ""
$ws.Range("P27").Value2 = ""
$ws.Range("Q27").Value2 = ""
$ws.Range("R27").Value2 = ""
$ws.Range("S27").Value2 = "Note di Stelle"
$ws.Range("T27").Value2 = ""
$ws.Range("U27").Value2 = ""
$ws.Range("V27").Value2 = $false
$ws.Range("W27").Value2 = 41123
$ws.Range("X27").Value2 = "https://www.comune.modena.it/novita/eventi/2022/note-di-stelle"
$ws.Range("Y27").Value2 = "44,64582"
$ws.Range("Z27").Value2 = "10,92572"
$ws.Range("AA27").Value2 = "POINT (10.92572 44.64582)"
